$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet's tab/title and update the sheet name reference in workbook
$ws.Name = "Through 2021-09-24"

# Update the label in A10 (shared string "September (through 09-23)" -> "September (through 09-24)")
$ws.Range("A10").Value = "September (through 09-24)"

# Update H8 value
$ws.Range("H8").Value = 151

# Update row 10 (September row) values
$ws.Range("C10").Value = 35
$ws.Range("D10").Value = 60
$ws.Range("E10").Value = 47
$ws.Range("F10").Value = 60
$ws.Range("G10").Value = 93
$ws.Range("H10").Value = 142

# Update row 11 (Total row) values
$ws.Range("C11").Value = 416
$ws.Range("D11").Value = 611
$ws.Range("E11").Value = 537
$ws.Range("F11").Value = 409
$ws.Range("G11").Value = 877
$ws.Range("H11").Value = 1212
